$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.6654078020343008
$ws.Cells.Item(2, 3).Value = -0.3131018370847687
$ws.Cells.Item(2, 4).Value = 0.3585108818504688
$ws.Cells.Item(2, 5).Value = 0.06866841829193859
$ws.Cells.Item(2, 6).Value = 0.370294988155365
$ws.Cells.Item(2, 7).Value = 1.238277912139893
$ws.Cells.Item(2, 8).Value = 0.8910083770751953
$ws.Cells.Item(2, 9).Value = 1.074856638908386

$ws.Cells.Item(3, 1).Value = "model_10_1_8"
$ws.Cells.Item(3, 2).Value = 0.6726587272939677
$ws.Cells.Item(3, 3).Value = -0.3748104032787467
$ws.Cells.Item(3, 4).Value = 0.3298698521879422
$ws.Cells.Item(3, 5).Value = 0.0257531521474389
$ws.Cells.Item(3, 6).Value = 0.362270325422287
$ws.Cells.Item(3, 7).Value = 1.296470046043396
$ws.Cells.Item(3, 8).Value = 0.9307898879051208
$ws.Cells.Item(3, 9).Value = 1.124385476112366

$ws.Cells.Item(4, 1).Value = "model_10_1_17"
$ws.Cells.Item(4, 2).Value = 0.6727424936888859
$ws.Cells.Item(4, 3).Value = -0.3645084027780809
$ws.Cells.Item(4, 4).Value = 0.3220506528079983
$ws.Cells.Item(4, 5).Value = 0.02578114564255063
$ws.Cells.Item(4, 6).Value = 0.3621776401996613
$ws.Cells.Item(4, 7).Value = 1.286755084991455
$ws.Cells.Item(4, 8).Value = 0.9416505098342896
$ws.Cells.Item(4, 9).Value = 1.124353051185608

$ws.Cells.Item(5, 2).Value = 0.6731329884640765
$ws.Cells.Item(5, 3).Value = -0.3571550579828167
$ws.Cells.Item(5, 4).Value = 0.3184464987679826
$ws.Cells.Item(5, 5).Value = 0.02692118404605293
$ws.Cells.Item(5, 6).Value = 0.3617455065250397
$ws.Cells.Item(5, 7).Value = 1.279820799827576
$ws.Cells.Item(5, 8).Value = 0.9466565847396851
$ws.Cells.Item(5, 9).Value = 1.123037338256836

$ws.Cells.Item(6, 2).Value = 0.6731329884640765
$ws.Cells.Item(6, 3).Value = -0.3571550579828167
$ws.Cells.Item(6, 4).Value = 0.3184464987679826
$ws.Cells.Item(6, 5).Value = 0.02692118404605293
$ws.Cells.Item(6, 6).Value = 0.3617455065250397
$ws.Cells.Item(6, 7).Value = 1.279820799827576
$ws.Cells.Item(6, 8).Value = 0.9466565847396851
$ws.Cells.Item(6, 9).Value = 1.123037338256836

$ws.Cells.Item(7, 1).Value = "model_10_1_9"
$ws.Cells.Item(7, 2).Value = 0.6731409431542253
$ws.Cells.Item(7, 3).Value = -0.3710434286522086
$ws.Cells.Item(7, 4).Value = 0.3296371343688055
$ws.Cells.Item(7, 5).Value = 0.0272508088946275
$ws.Cells.Item(7, 6).Value = 0.3617366850376129
$ws.Cells.Item(7, 7).Value = 1.292917728424072
$ws.Cells.Item(7, 8).Value = 0.931113064289093
$ws.Cells.Item(7, 9).Value = 1.12265682220459

$ws.Cells.Item(8, 1).Value = "model_10_1_22"
$ws.Cells.Item(8, 2).Value = 0.6731504177283478
$ws.Cells.Item(8, 3).Value = -0.357095672582511
$ws.Cells.Item(8, 4).Value = 0.3184986020383779
$ws.Cells.Item(8, 5).Value = 0.0269764204414833
$ws.Cells.Item(8, 6).Value = 0.3617261350154877
$ws.Cells.Item(8, 7).Value = 1.279764652252197
$ws.Cells.Item(8, 8).Value = 0.946584165096283
$ws.Cells.Item(8, 9).Value = 1.122973561286926

$ws.Cells.Item(9, 1).Value = "model_10_1_21"
$ws.Cells.Item(9, 2).Value = 0.6732500778038092
$ws.Cells.Item(9, 3).Value = -0.3568397586847933
$ws.Cells.Item(9, 4).Value = 0.3188706310995222
$ws.Cells.Item(9, 5).Value = 0.02729768170598912
$ws.Cells.Item(9, 6).Value = 0.3616158664226532
$ws.Cells.Item(9, 7).Value = 1.279523372650146
$ws.Cells.Item(9, 8).Value = 0.9460674524307251
$ws.Cells.Item(9, 9).Value = 1.122602939605713

$ws.Cells.Item(10, 1).Value = "model_10_1_4"
$ws.Cells.Item(10, 2).Value = 0.6733657827974664
$ws.Cells.Item(10, 3).Value = -0.408839100924822
$ws.Cells.Item(10, 4).Value = 0.3647972904289529
$ws.Cells.Item(10, 5).Value = 0.03081439831803179
$ws.Cells.Item(10, 6).Value = 0.3614878356456757
$ws.Cells.Item(10, 7).Value = 1.328559637069702
$ws.Cells.Item(10, 8).Value = 0.882276713848114
$ws.Cells.Item(10, 9).Value = 1.118544220924377

$ws.Cells.Item(11, 2).Value = 0.673532426651683
$ws.Cells.Item(11, 3).Value = -0.3829196245382855
$ws.Cells.Item(11, 4).Value = 0.3439723656436888
$ws.Cells.Item(11, 5).Value = 0.03023239865027927
$ws.Cells.Item(11, 6).Value = 0.3613033890724182
$ws.Cells.Item(11, 7).Value = 1.304117202758789
$ws.Cells.Item(11, 8).Value = 0.9112019538879395
$ws.Cells.Item(11, 9).Value = 1.119215965270996

$ws.Cells.Item(12, 1).Value = "model_10_1_20"
$ws.Cells.Item(12, 2).Value = 0.6735494200212875
$ws.Cells.Item(12, 3).Value = -0.3550319634784489
$ws.Cells.Item(12, 4).Value = 0.3191726796548366
$ws.Cells.Item(12, 5).Value = 0.02825070305008059
$ws.Cells.Item(12, 6).Value = 0.3612845838069916
$ws.Cells.Item(12, 7).Value = 1.277818560600281
$ws.Cells.Item(12, 8).Value = 0.9456478953361511
$ws.Cells.Item(12, 9).Value = 1.121502876281738

$ws.Cells.Item(13, 1).Value = "model_10_1_15"
$ws.Cells.Item(13, 2).Value = 0.673687612309118
$ws.Cells.Item(13, 3).Value = -0.3595737632941307
$ws.Cells.Item(13, 4).Value = 0.3235884425490637
$ws.Cells.Item(13, 5).Value = 0.02878693929150922
$ws.Cells.Item(13, 6).Value = 0.3611316680908203
$ws.Cells.Item(13, 7).Value = 1.282101631164551
$ws.Cells.Item(13, 8).Value = 0.9395145773887634
$ws.Cells.Item(13, 9).Value = 1.12088406085968

$ws.Cells.Item(14, 2).Value = 0.6737165394695905
$ws.Cells.Item(14, 3).Value = -0.3588336186717982
$ws.Cells.Item(14, 4).Value = 0.3231735537513322
$ws.Cells.Item(14, 5).Value = 0.02887196263485592
$ws.Cells.Item(14, 6).Value = 0.3610996603965759
$ws.Cells.Item(14, 7).Value = 1.281403660774231
$ws.Cells.Item(14, 8).Value = 0.9400908350944519
$ws.Cells.Item(14, 9).Value = 1.12078595161438

$ws.Cells.Item(15, 1).Value = "model_10_1_10"
$ws.Cells.Item(15, 2).Value = 0.6739146618586163
$ws.Cells.Item(15, 3).Value = -0.3601612090442967
$ws.Cells.Item(15, 4).Value = 0.3254017354607431
$ws.Cells.Item(15, 5).Value = 0.02955967849416208
$ws.Cells.Item(15, 6).Value = 0.3608803749084473
$ws.Cells.Item(15, 7).Value = 1.282655477523804
$ws.Cells.Item(15, 8).Value = 0.9369959235191345
$ws.Cells.Item(15, 9).Value = 1.119992256164551

$ws.Cells.Item(16, 1).Value = "model_10_1_18"
$ws.Cells.Item(16, 2).Value = 0.6739924046657439
$ws.Cells.Item(16, 3).Value = -0.3541650921807626
$ws.Cells.Item(16, 4).Value = 0.3210473704135725
$ws.Cells.Item(16, 5).Value = 0.02968734599361322
$ws.Cells.Item(16, 6).Value = 0.360794335603714
$ws.Cells.Item(16, 7).Value = 1.277001142501831
$ws.Cells.Item(16, 8).Value = 0.9430440068244934
$ws.Cells.Item(16, 9).Value = 1.119845032691956

$ws.Cells.Item(17, 1).Value = "model_10_1_11"
$ws.Cells.Item(17, 2).Value = 0.6740106434606149
$ws.Cells.Item(17, 3).Value = -0.3592224748611699
$ws.Cells.Item(17, 4).Value = 0.325212159957798
$ws.Cells.Item(17, 5).Value = 0.0298583570729779
$ws.Cells.Item(17, 6).Value = 0.3607741892337799
$ws.Cells.Item(17, 7).Value = 1.281770348548889
$ws.Cells.Item(17, 8).Value = 0.9372592568397522
$ws.Cells.Item(17, 9).Value = 1.119647622108459

$ws.Cells.Item(18, 2).Value = 0.6745474870293133
$ws.Cells.Item(18, 3).Value = -0.3538835773135671
$ws.Cells.Item(18, 4).Value = 0.3241258585839981
$ws.Cells.Item(18, 5).Value = 0.03155272047093471
$ws.Cells.Item(18, 6).Value = 0.3601800501346588
$ws.Cells.Item(18, 7).Value = 1.276735782623291
$ws.Cells.Item(18, 8).Value = 0.9387680888175964
$ws.Cells.Item(18, 9).Value = 1.117691993713379

$ws.Cells.Item(19, 1).Value = "model_10_1_19"
$ws.Cells.Item(19, 2).Value = 0.6748358163897757
$ws.Cells.Item(19, 3).Value = -0.3476449633534389
$ws.Cells.Item(19, 4).Value = 0.3207594424570495
$ws.Cells.Item(19, 5).Value = 0.03234475863462904
$ws.Cells.Item(19, 6).Value = 0.3598609268665314
$ws.Cells.Item(19, 7).Value = 1.270852565765381
$ws.Cells.Item(19, 8).Value = 0.9434439539909363
$ws.Cells.Item(19, 9).Value = 1.116778016090393

$ws.Cells.Item(20, 1).Value = "model_10_1_13"
$ws.Cells.Item(20, 2).Value = 0.6748673028492558
$ws.Cells.Item(20, 3).Value = -0.3517407435665285
$ws.Cells.Item(20, 4).Value = 0.3242274567428035
$ws.Cells.Item(20, 5).Value = 0.03253697434324487
$ws.Cells.Item(20, 6).Value = 0.3598260879516602
$ws.Cells.Item(20, 7).Value = 1.274714946746826
$ws.Cells.Item(20, 8).Value = 0.9386269450187683
$ws.Cells.Item(20, 9).Value = 1.116556286811829

$ws.Cells.Item(21, 1).Value = "model_10_1_14"
$ws.Cells.Item(21, 2).Value = 0.6750390411828946
$ws.Cells.Item(21, 3).Value = -0.3500642146970379
$ws.Cells.Item(21, 4).Value = 0.3238914944442692
$ws.Cells.Item(21, 5).Value = 0.03307228979176946
$ws.Cells.Item(21, 6).Value = 0.3596360385417938
$ws.Cells.Item(21, 7).Value = 1.273133873939514
$ws.Cells.Item(21, 8).Value = 0.9390936493873596
$ws.Cells.Item(21, 9).Value = 1.115938305854797

$ws.Cells.Item(22, 1).Value = "model_10_1_7"
$ws.Cells.Item(22, 2).Value = 0.6754390636708122
$ws.Cells.Item(22, 3).Value = -0.3649794492261731
$ws.Cells.Item(22, 4).Value = 0.3383637711592755
$ws.Cells.Item(22, 5).Value = 0.03481631750169156
$ws.Cells.Item(22, 6).Value = 0.3591933250427246
$ws.Cells.Item(22, 7).Value = 1.287199378013611
$ws.Cells.Item(22, 8).Value = 0.9189921021461487
$ws.Cells.Item(22, 9).Value = 1.113925576210022

$ws.Cells.Item(23, 1).Value = "model_10_1_6"
$ws.Cells.Item(23, 2).Value = 0.6755164568056855
$ws.Cells.Item(23, 3).Value = -0.3656297570528668
$ws.Cells.Item(23, 4).Value = 0.3394384174549508
$ws.Cells.Item(23, 5).Value = 0.03514396526578034
$ws.Cells.Item(23, 6).Value = 0.3591076731681824
$ws.Cells.Item(23, 7).Value = 1.28781259059906
$ws.Cells.Item(23, 8).Value = 0.9174994230270386
$ws.Cells.Item(23, 9).Value = 1.113547444343567

$ws.Cells.Item(24, 1).Value = "model_10_1_3"
$ws.Cells.Item(24, 2).Value = 0.6859054242072302
$ws.Cells.Item(24, 3).Value = -0.2837763372652367
$ws.Cells.Item(24, 4).Value = 0.3611801751803887
$ws.Cells.Item(24, 5).Value = 0.08286540790518027
$ws.Cells.Item(24, 6).Value = 0.3476101458072662
$ws.Cells.Item(24, 7).Value = 1.210623264312744
$ws.Cells.Item(24, 8).Value = 0.8873007893562317
$ws.Cells.Item(24, 9).Value = 1.05847179889679

$ws.Cells.Item(25, 2).Value = 0.7267443074194322
$ws.Cells.Item(25, 3).Value = 0.04213534430468391
$ws.Cells.Item(25, 4).Value = 0.5526929139258074
$ws.Cells.Item(25, 5).Value = 0.3323125230842148
$ws.Cells.Item(25, 6).Value = 0.3024135231971741
$ws.Cells.Item(25, 7).Value = 0.9032828807830811
$ws.Cells.Item(25, 8).Value = 0.6212955713272095
$ws.Cells.Item(25, 9).Value = 0.7705829739570618

$ws.Cells.Item(26, 2).Value = 0.7463476027817892
$ws.Cells.Item(26, 3).Value = 0.4785403134894666
$ws.Cells.Item(26, 4).Value = 0.6323103106306069
$ws.Cells.Item(26, 5).Value = 0.5661844059404084
$ws.Cells.Item(26, 6).Value = 0.2807184755802155
$ws.Cells.Item(26, 7).Value = 0.4917455315589905
$ws.Cells.Item(26, 8).Value = 0.5107095241546631
$ws.Cells.Item(26, 9).Value = 0.5006697773933411
